# Update the "取得日時" (acquisition timestamp) column for the job-listing
# rows on the "ランサーズ" sheet: append a refreshed timestamp of
# 2026-01-12 12:55:58 to rows 2-8 (column A), replacing the previous
# 2026-01-12 12:42:28 value, mirroring a re-scrape/append run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-12 12:55:58"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
